$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.743.51'
$ws.Range('E2').Value = '  -1.77%  '
$ws.Range('D3').Value = '3.484.51'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.14'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.00%  '
$ws.Range('D7').Value = '3.483.70'
$ws.Range('E7').Value = '  -0.96%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -2.02%  '
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.16'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.380'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('D13').Value = '4.072.96'
$ws.Range('E13').Value = '  -0.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.32'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('E15').Value = '  +1.39%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000176'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.39%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.468.66'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').Value = '63.822.32'
$ws.Range('E18').Value = '  -1.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.19'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.64'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '383.31'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.03%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').Value = '3.621.65'
$ws.Range('E24').Value = '  -1.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.93%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  +1.79%  '
$ws.Range('E28').Value = '  -0.98%  '
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.28'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('E32').Value = '  -2.30%  '
$ws.Range('D33').Value = '3.491.21'
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.46'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.14%  '
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('E37').Value = '  +4.53%  '
$ws.Range('E38').Value = '  +2.11%  '
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '160.05'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0795'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.47'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.39%  '
$ws.Range('E43').Value = '  -1.46%  '
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.76'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.30%  '
$ws.Range('E46').Value = '  -2.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.39'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('E48').Value = '  -1.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.84'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('D50').Value = '2.417.08'
$ws.Range('E50').Value = '  +1.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.901'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.71%  '
